$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.616.23'
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").Value = '2.945.25'
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '364.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.73%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.75%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -3.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.140'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0844'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.65%  '

$ws.Range("D14").Value = '3.410.95'
$ws.Range("E14").Value = '  +0.98%  '

$ws.Range("E15").Value = '  -4.44%  '

$ws.Range("D16").Value = '2.955.98'
$ws.Range("E16").Value = '  +1.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.980'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("D18").Value = '51.537.53'
$ws.Range("E18").Value = '  -0.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("E21").Value = '  -4.58%  '

$ws.Range("D22").Value = '0.0₃0955'
$ws.Range("E22").Value = '  -2.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.58%  '

$ws.Range("E25").Value = '  -3.22%  '

$ws.Range("E26").Value = '  -5.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.111'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.79%  '

$ws.Range("E33").Value = '  +5.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0429'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.80%  '

$ws.Range("E37").Value = '  +0.22%  '

$ws.Range("E38").Value = '  +4.73%  '

$ws.Range("E39").Value = '  +0.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.12%  '

$ws.Range("E41").Value = '  -4.76%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.115'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.64%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '120.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.13%  '

$ws.Range("E45").Value = '  -1.18%  '

$ws.Range("D46").Value = '2.101.33'
$ws.Range("E46").Value = '  -0.81%  '

$ws.Range("E47").Value = '  -5.90%  '

$ws.Range("E48").Value = '  -7.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.240'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.33%  '

$ws.Range("E50").Value = '  -5.91%  '

$ws.Range("E51").Value = '  -2.59%  '

